$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 'BTC'
$ws.Cells.Item(2, 3).Value = 'Bitcoin'
$ws.Cells.Item(2, 4).Value = 30691
$ws.Cells.Item(2, 5).Value = 595759907750
$ws.Cells.Item(2, 6).Value = 12975786163
$ws.Cells.Item(2, 7).Value = -0.03812
$ws.Cells.Item(3, 2).Value = 'ETH'
$ws.Cells.Item(3, 3).Value = 'Ethereum'
$ws.Cells.Item(3, 4).Value = 1916.11
$ws.Cells.Item(3, 5).Value = 230285731243
$ws.Cells.Item(3, 6).Value = 7280425873
$ws.Cells.Item(3, 7).Value = 1.02868
$ws.Cells.Item(4, 2).Value = 'USDT'
$ws.Cells.Item(4, 3).Value = 'Tether'
$ws.Cells.Item(4, 4).Value = 1
$ws.Cells.Item(4, 5).Value = 83130016547
$ws.Cells.Item(4, 6).Value = 19416891910
$ws.Cells.Item(4, 7).Value = -0.05126
$ws.Cells.Item(5, 2).Value = 'BNB'
$ws.Cells.Item(5, 3).Value = 'BNB'
$ws.Cells.Item(5, 4).Value = 239.95
$ws.Cells.Item(5, 5).Value = 37397549306
$ws.Cells.Item(5, 6).Value = 766763641
$ws.Cells.Item(5, 7).Value = -2.16308
$ws.Cells.Item(6, 2).Value = 'USDC'
$ws.Cells.Item(6, 3).Value = 'USD Coin'
$ws.Cells.Item(6, 4).Value = 1
$ws.Cells.Item(6, 5).Value = 28436790925
$ws.Cells.Item(6, 6).Value = 3155463040
$ws.Cells.Item(6, 7).Value = -0.00446
$ws.Cells.Item(7, 2).Value = 'XRP'
$ws.Cells.Item(7, 3).Value = 'XRP'
$ws.Cells.Item(7, 4).Value = 0.491939
$ws.Cells.Item(7, 5).Value = 25713728479
$ws.Cells.Item(7, 6).Value = 542038095
$ws.Cells.Item(7, 7).Value = -0.1495
$ws.Cells.Item(8, 2).Value = 'STETH'
$ws.Cells.Item(8, 3).Value = 'Lido Staked Ether'
$ws.Cells.Item(8, 4).Value = 1912.08
$ws.Cells.Item(8, 5).Value = 14141089162
$ws.Cells.Item(8, 6).Value = 2000483
$ws.Cells.Item(8, 7).Value = 0.81103
$ws.Cells.Item(9, 2).Value = 'ADA'
$ws.Cells.Item(9, 3).Value = 'Cardano'
$ws.Cells.Item(9, 4).Value = 0.297543
$ws.Cells.Item(9, 5).Value = 10433002658
$ws.Cells.Item(9, 6).Value = 210087347
$ws.Cells.Item(9, 7).Value = 0.56745
$ws.Cells.Item(10, 2).Value = 'DOGE'
$ws.Cells.Item(10, 3).Value = 'Dogecoin'
$ws.Cells.Item(10, 4).Value = 0.06768399999999999
$ws.Cells.Item(10, 5).Value = 9468276785
$ws.Cells.Item(10, 6).Value = 243364562
$ws.Cells.Item(10, 7).Value = -0.43341
$ws.Cells.Item(11, 2).Value = 'SOL'
$ws.Cells.Item(11, 3).Value = 'Solana'
$ws.Cells.Item(11, 4).Value = 17.21
$ws.Cells.Item(11, 5).Value = 6887534070
$ws.Cells.Item(11, 6).Value = 248616422
$ws.Cells.Item(11, 7).Value = -0.61104
$ws.Cells.Item(12, 2).Value = 'TRX'
$ws.Cells.Item(12, 3).Value = 'TRON'
$ws.Cells.Item(12, 4).Value = 0.073522
$ws.Cells.Item(12, 5).Value = 6603860079
$ws.Cells.Item(12, 6).Value = 246993465
$ws.Cells.Item(12, 7).Value = 1.38104
$ws.Cells.Item(13, 2).Value = 'LTC'
$ws.Cells.Item(13, 3).Value = 'Litecoin'
$ws.Cells.Item(13, 4).Value = 88.84
$ws.Cells.Item(13, 5).Value = 6507691635
$ws.Cells.Item(13, 6).Value = 587917523
$ws.Cells.Item(13, 7).Value = -2.40692
$ws.Cells.Item(14, 2).Value = 'DOT'
$ws.Cells.Item(14, 3).Value = 'Polkadot'
$ws.Cells.Item(14, 4).Value = 5.19
$ws.Cells.Item(14, 5).Value = 6456630286
$ws.Cells.Item(14, 6).Value = 127399249
$ws.Cells.Item(14, 7).Value = 2.67859
$ws.Cells.Item(15, 2).Value = 'MATIC'
$ws.Cells.Item(15, 3).Value = 'Polygon'
$ws.Cells.Item(15, 4).Value = 0.672648
$ws.Cells.Item(15, 5).Value = 6251692478
$ws.Cells.Item(15, 6).Value = 228944836
$ws.Cells.Item(15, 7).Value = -1.10279
$ws.Cells.Item(16, 2).Value = 'WBTC'
$ws.Cells.Item(16, 3).Value = 'Wrapped Bitcoin'
$ws.Cells.Item(16, 4).Value = 30646
$ws.Cells.Item(16, 5).Value = 4818705018
$ws.Cells.Item(16, 6).Value = 151713562
$ws.Cells.Item(16, 7).Value = -0.27449
$ws.Cells.Item(17, 2).Value = 'SHIB'
$ws.Cells.Item(17, 3).Value = 'Shiba Inu'
$ws.Cells.Item(17, 4).Value = 0.00000795
$ws.Cells.Item(17, 5).Value = 4681027830
$ws.Cells.Item(17, 6).Value = 116390454
$ws.Cells.Item(17, 7).Value = -1.15743
$ws.Cells.Item(18, 2).Value = 'AVAX'
$ws.Cells.Item(18, 3).Value = 'Avalanche'
$ws.Cells.Item(18, 4).Value = 13.55
$ws.Cells.Item(18, 5).Value = 4677112416
$ws.Cells.Item(18, 6).Value = 139895444
$ws.Cells.Item(18, 7).Value = 2.69807
$ws.Cells.Item(19, 2).Value = 'DAI'
$ws.Cells.Item(19, 3).Value = 'Dai'
$ws.Cells.Item(19, 4).Value = 1
$ws.Cells.Item(19, 5).Value = 4450286885
$ws.Cells.Item(19, 6).Value = 87921293
$ws.Cells.Item(19, 7).Value = -0.07027
$ws.Cells.Item(20, 2).Value = 'BUSD'
$ws.Cells.Item(20, 3).Value = 'Binance USD'
$ws.Cells.Item(20, 4).Value = 1
$ws.Cells.Item(20, 5).Value = 4244868393
$ws.Cells.Item(20, 6).Value = 2106858046
$ws.Cells.Item(20, 7).Value = -0.05094
$ws.Cells.Item(21, 2).Value = 'UNI'
$ws.Cells.Item(21, 3).Value = 'Uniswap'
$ws.Cells.Item(21, 4).Value = 5.38
$ws.Cells.Item(21, 5).Value = 4030766043
$ws.Cells.Item(21, 6).Value = 109292465
$ws.Cells.Item(21, 7).Value = 11.56179
$ws.Cells.Item(22, 2).Value = 'BCH'
$ws.Cells.Item(22, 3).Value = 'Bitcoin Cash'
$ws.Cells.Item(22, 4).Value = 200.98
$ws.Cells.Item(22, 5).Value = 3905237719
$ws.Cells.Item(22, 6).Value = 1363922665
$ws.Cells.Item(22, 7).Value = 3.75765
$ws.Cells.Item(23, 2).Value = 'LEO'
$ws.Cells.Item(23, 3).Value = 'LEO Token'
$ws.Cells.Item(23, 4).Value = 3.84
$ws.Cells.Item(23, 5).Value = 3571855159
$ws.Cells.Item(23, 6).Value = 713801
$ws.Cells.Item(23, 7).Value = -1.00332
$ws.Cells.Item(24, 2).Value = 'LINK'
$ws.Cells.Item(24, 3).Value = 'Chainlink'
$ws.Cells.Item(24, 4).Value = 6.3
$ws.Cells.Item(24, 5).Value = 3262044470
$ws.Cells.Item(24, 6).Value = 166878918
$ws.Cells.Item(24, 7).Value = 2.4205
$ws.Cells.Item(25, 2).Value = 'TUSD'
$ws.Cells.Item(25, 3).Value = 'TrueUSD'
$ws.Cells.Item(25, 4).Value = 0.99897
$ws.Cells.Item(25, 5).Value = 3135245545
$ws.Cells.Item(25, 6).Value = 189437775
$ws.Cells.Item(25, 7).Value = -0.10681
$ws.Cells.Item(26, 2).Value = 'XMR'
$ws.Cells.Item(26, 3).Value = 'Monero'
$ws.Cells.Item(26, 4).Value = 165.52
$ws.Cells.Item(26, 5).Value = 3011967273
$ws.Cells.Item(26, 6).Value = 95791616
$ws.Cells.Item(26, 7).Value = 6.16086
$ws.Cells.Item(27, 2).Value = 'ATOM'
$ws.Cells.Item(27, 3).Value = 'Cosmos Hub'
$ws.Cells.Item(27, 4).Value = 9.66
$ws.Cells.Item(27, 5).Value = 2822832048
$ws.Cells.Item(27, 6).Value = 164204223
$ws.Cells.Item(27, 7).Value = 2.52374
$ws.Cells.Item(28, 2).Value = 'OKB'
$ws.Cells.Item(28, 3).Value = 'OKB'
$ws.Cells.Item(28, 4).Value = 44.96
$ws.Cells.Item(28, 5).Value = 2705385551
$ws.Cells.Item(28, 6).Value = 4718709
$ws.Cells.Item(28, 7).Value = -0.34113
$ws.Cells.Item(29, 2).Value = 'ETC'
$ws.Cells.Item(29, 3).Value = 'Ethereum Classic'
$ws.Cells.Item(29, 4).Value = 18.9
$ws.Cells.Item(29, 5).Value = 2672922535
$ws.Cells.Item(29, 6).Value = 243135797
$ws.Cells.Item(29, 7).Value = -1.32338
$ws.Cells.Item(30, 2).Value = 'XLM'
$ws.Cells.Item(30, 3).Value = 'Stellar'
$ws.Cells.Item(30, 4).Value = 0.091516
$ws.Cells.Item(30, 5).Value = 2466137755
$ws.Cells.Item(30, 6).Value = 42374487
$ws.Cells.Item(30, 7).Value = 0.42965
$ws.Cells.Item(31, 2).Value = 'TON'
$ws.Cells.Item(31, 3).Value = 'Toncoin'
$ws.Cells.Item(31, 4).Value = 1.48
$ws.Cells.Item(31, 5).Value = 2175995111
$ws.Cells.Item(31, 6).Value = 8993026
$ws.Cells.Item(31, 7).Value = 5.49221
$ws.Cells.Item(32, 2).Value = 'ICP'
$ws.Cells.Item(32, 3).Value = 'Internet Computer'
$ws.Cells.Item(32, 4).Value = 4.36
$ws.Cells.Item(32, 5).Value = 1910550258
$ws.Cells.Item(32, 6).Value = 14987941
$ws.Cells.Item(32, 7).Value = 0.67437
$ws.Cells.Item(33, 2).Value = 'FIL'
$ws.Cells.Item(33, 3).Value = 'Filecoin'
$ws.Cells.Item(33, 4).Value = 4.08
$ws.Cells.Item(33, 5).Value = 1757764243
$ws.Cells.Item(33, 6).Value = 91776724
$ws.Cells.Item(33, 7).Value = 1.22967
$ws.Cells.Item(34, 2).Value = 'LDO'
$ws.Cells.Item(34, 3).Value = 'Lido DAO'
$ws.Cells.Item(34, 4).Value = 1.97
$ws.Cells.Item(34, 5).Value = 1728499971
$ws.Cells.Item(34, 6).Value = 38411138
$ws.Cells.Item(34, 7).Value = 3.06545
$ws.Cells.Item(35, 2).Value = 'HBAR'
$ws.Cells.Item(35, 3).Value = 'Hedera'
$ws.Cells.Item(35, 4).Value = 0.053137
$ws.Cells.Item(35, 5).Value = 1685966481
$ws.Cells.Item(35, 6).Value = 40233097
$ws.Cells.Item(35, 7).Value = 1.61432
$ws.Cells.Item(36, 2).Value = 'APT'
$ws.Cells.Item(36, 3).Value = 'Aptos'
$ws.Cells.Item(36, 4).Value = 7.65
$ws.Cells.Item(36, 5).Value = 1589465079
$ws.Cells.Item(36, 6).Value = 88098777
$ws.Cells.Item(36, 7).Value = 0.17505
$ws.Cells.Item(37, 2).Value = 'QNT'
$ws.Cells.Item(37, 3).Value = 'Quant'
$ws.Cells.Item(37, 4).Value = 107.19
$ws.Cells.Item(37, 5).Value = 1556857575
$ws.Cells.Item(37, 6).Value = 15915919
$ws.Cells.Item(37, 7).Value = 1.35777
$ws.Cells.Item(38, 2).Value = 'CRO'
$ws.Cells.Item(38, 3).Value = 'Cronos'
$ws.Cells.Item(38, 4).Value = 0.059079
$ws.Cells.Item(38, 5).Value = 1544761405
$ws.Cells.Item(38, 6).Value = 9708126
$ws.Cells.Item(38, 7).Value = 0.69263
$ws.Cells.Item(39, 2).Value = 'ARB'
$ws.Cells.Item(39, 3).Value = 'Arbitrum'
$ws.Cells.Item(39, 4).Value = 1.12
$ws.Cells.Item(39, 5).Value = 1428066049
$ws.Cells.Item(39, 6).Value = 189993603
$ws.Cells.Item(39, 7).Value = 0.93971
$ws.Cells.Item(40, 2).Value = 'NEAR'
$ws.Cells.Item(40, 3).Value = 'NEAR Protocol'
$ws.Cells.Item(40, 4).Value = 1.44
$ws.Cells.Item(40, 5).Value = 1336181561
$ws.Cells.Item(40, 6).Value = 46751318
$ws.Cells.Item(40, 7).Value = 1.02952
$ws.Cells.Item(41, 2).Value = 'VET'
$ws.Cells.Item(41, 3).Value = 'VeChain'
$ws.Cells.Item(41, 4).Value = 0.0183611
$ws.Cells.Item(41, 5).Value = 1335366878
$ws.Cells.Item(41, 6).Value = 51254064
$ws.Cells.Item(41, 7).Value = -0.07312
$ws.Cells.Item(42, 2).Value = 'AAVE'
$ws.Cells.Item(42, 3).Value = 'Aave'
$ws.Cells.Item(42, 4).Value = 76.27
$ws.Cells.Item(42, 5).Value = 1101005817
$ws.Cells.Item(42, 6).Value = 292426031
$ws.Cells.Item(42, 7).Value = 32.00964
$ws.Cells.Item(43, 2).Value = 'GRT'
$ws.Cells.Item(43, 3).Value = 'The Graph'
$ws.Cells.Item(43, 4).Value = 0.114258
$ws.Cells.Item(43, 5).Value = 1033879778
$ws.Cells.Item(43, 6).Value = 41695558
$ws.Cells.Item(43, 7).Value = 2.5478
$ws.Cells.Item(44, 2).Value = 'STX'
$ws.Cells.Item(44, 3).Value = 'Stacks'
$ws.Cells.Item(44, 4).Value = 0.732878
$ws.Cells.Item(44, 5).Value = 1019377709
$ws.Cells.Item(44, 6).Value = 34354317
$ws.Cells.Item(44, 7).Value = -2.65856
$ws.Cells.Item(45, 2).Value = 'ALGO'
$ws.Cells.Item(45, 3).Value = 'Algorand'
$ws.Cells.Item(45, 4).Value = 0.138993
$ws.Cells.Item(45, 5).Value = 1007806744
$ws.Cells.Item(45, 6).Value = 31973157
$ws.Cells.Item(45, 7).Value = 3.2671
$ws.Cells.Item(46, 2).Value = 'FRAX'
$ws.Cells.Item(46, 3).Value = 'Frax'
$ws.Cells.Item(46, 4).Value = 0.999351
$ws.Cells.Item(46, 5).Value = 1002415962
$ws.Cells.Item(46, 6).Value = 7947944
$ws.Cells.Item(46, 7).Value = -0.11412
$ws.Cells.Item(47, 2).Value = 'USDP'
$ws.Cells.Item(47, 3).Value = 'Pax Dollar'
$ws.Cells.Item(47, 4).Value = 0.99935
$ws.Cells.Item(47, 5).Value = 999720595
$ws.Cells.Item(47, 6).Value = 1221012
$ws.Cells.Item(47, 7).Value = -0.04413
$ws.Cells.Item(48, 2).Value = 'RETH'
$ws.Cells.Item(48, 3).Value = 'Rocket Pool ETH'
$ws.Cells.Item(48, 4).Value = 2060.91
$ws.Cells.Item(48, 5).Value = 945129497
$ws.Cells.Item(48, 6).Value = 2297779
$ws.Cells.Item(48, 7).Value = 1.10236
$ws.Cells.Item(49, 2).Value = 'EGLD'
$ws.Cells.Item(49, 3).Value = 'MultiversX'
$ws.Cells.Item(49, 4).Value = 35.99
$ws.Cells.Item(49, 5).Value = 918222507
$ws.Cells.Item(49, 6).Value = 27069294
$ws.Cells.Item(49, 7).Value = 6.69028
$ws.Cells.Item(50, 2).Value = 'FTM'
$ws.Cells.Item(50, 3).Value = 'Fantom'
$ws.Cells.Item(50, 4).Value = 0.326793
$ws.Cells.Item(50, 5).Value = 916300078
$ws.Cells.Item(50, 6).Value = 122537258
$ws.Cells.Item(50, 7).Value = 2.01302
$ws.Cells.Item(51, 2).Value = 'OP'
$ws.Cells.Item(51, 3).Value = 'Optimism'
$ws.Cells.Item(51, 4).Value = 1.36
$ws.Cells.Item(51, 5).Value = 879222145
$ws.Cells.Item(51, 6).Value = 103883164
$ws.Cells.Item(51, 7).Value = 1.8078
